$d = $word.ActiveDocument

# --- 1. New "Abstract Title" paragraph style, inserted (conceptually) right
#        before the existing "Abstract" style. ------------------------------
$abstractTitle = $d.Styles.Add("AbstractTitle", 1)
$abstractTitle.NameLocal = "Abstract Title"
$abstractTitle.BaseStyle = "Normal"
$abstractTitle.NextParagraphStyle = "Abstract"
$abstractTitle.QuickStyle = $true

$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.Alignment = 1
$abstractTitle.ParagraphFormat.SpaceAfter = 0
$abstractTitle.ParagraphFormat.SpaceBefore = 15

$abstractTitle.Font.Size = 10
$abstractTitle.Font.SizeBi = 10
$abstractTitle.Font.Bold = $true
$abstractTitle.Font.Color = 9067060

# --- 2. Existing "Abstract" style: reduce the space-before from 15pt (300)
#        to 5pt (100); space-after stays at 15pt (300). ---------------------
$abstract = $d.Styles.Add("Abstract", 1)
$abstract.ParagraphFormat.SpaceBefore = 5

# --- 3. New "Footnote Block Text" paragraph style, based on / followed by
#        "Footnote Text" (mirrors the existing "Block Text" style). --------
$footnoteBlockText = $d.Styles.Add("FootnoteBlockText", 1)
$footnoteBlockText.NameLocal = "Footnote Block Text"
$footnoteBlockText.BaseStyle = "Footnote Text"
$footnoteBlockText.NextParagraphStyle = "Footnote Text"
$footnoteBlockText.QuickStyle = $true
$footnoteBlockText.Priority = 9
$footnoteBlockText.UnhideWhenUsed = $true

$footnoteBlockText.ParagraphFormat.SpaceBefore = 5
$footnoteBlockText.ParagraphFormat.SpaceAfter = 5
$footnoteBlockText.ParagraphFormat.FirstLineIndent = 0
$footnoteBlockText.ParagraphFormat.LeftIndent = 24
$footnoteBlockText.ParagraphFormat.RightIndent = 24

Write-Host "Styles updated"
